# Apply the "Completed update to 1st iteration coding and data" edit.
# The study previously labelled "Greene2022#" in rows 12 and 13 is split
# into two distinctly-named studies: "Greene2022a" (row 12) and
# "Greene2022b" (row 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A12").Value = "Greene2022a"
$ws.Range("A13").Value = "Greene2022b"

# Restore the active selection to the cell it ended up on after editing.
$ws.Range("H11").Select()

$wb.Save()
